$wb = $excel.ActiveWorkbook

# --- Sheet1: rename to PB05_Login, add a new row of credentials ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "PB05_Login"

$ws1.Range("A5").Value = "admin@gmail.com"
$ws1.Range("B5").Value = "Password1%"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "mailto:admin@gmail.com")
$ws1.Range("A5").Style = "Hyperlink"

# --- New Sheet2: PB01_CreateStaffAccount, positioned right after sheet1 ---
# Duplicate sheet1 so the new sheet inherits matching row-height / column
# formatting, then strip it down and fill in the real content.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "PB01_CreateStaffAccount"

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = "StaffAccountTest@gmail.com"
$ws2.Range("A1").Value = "Email Created"
$ws2.Range("B1").Value = "Password Created"
$ws2.Range("B2").Value = 1234567890
$ws2.Range("A3:B5").Clear()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:StaffAccountTest@gmail.com")
$ws2.Range("A2").Style = "Hyperlink"

$ws2.Columns.Item(1).ColumnWidth = 25.67
$ws2.Columns.Item(2).ColumnWidth = 14.83

# --- Selections on each sheet (stored per-sheetView even when inactive) ---
$ws1.Range("C29").Select() | Out-Null
$ws2.Range("E20").Select() | Out-Null

# --- Make the new sheet the active tab ---
$ws2.Activate() | Out-Null

Write-Output "done"
